# "updated blank template spacing"
#
# The logo/title graphic on slide 1 (Oval 8, Oval 9, Freeform 10 and
# TextBox 11 - all children of the "Group 2" group shape) is shifted
# upward inside its (unchanged) group bounding box, tightening the
# vertical spacing used by the blank template.
#
# PowerPoint's Shape.Top/Left/Width/Height are expressed in points
# (1 pt = 12700 EMU) and are stored internally as single-precision
# floats, so the literal point values below are chosen to round-trip
# to the exact target EMU offsets from the authoritative OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$group = $s.Shapes.Item("Group 2")

$oval8      = $group.GroupItems.Item("Oval 8")
$oval9      = $group.GroupItems.Item("Oval 9")
$freeform10 = $group.GroupItems.Item("Freeform 10")
$textbox11  = $group.GroupItems.Item("TextBox 11")

# y: 1173613 EMU -> 440835 EMU (left/width/height untouched)
$oval8.Top = 34.71141732283465

# y: 1518414 EMU -> 785636 EMU
$oval9.Top = 61.861102362204726

# y: 3583315 EMU -> 2850537 EMU
$freeform10.Top = 224.45173228346457

# y: 1735841 EMU -> 1311279 EMU
$textbox11.Top = 103.25031666062992
